$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = 158.97082944738966
$ws.Range("C5").Value = 741.8950961577325
$ws.Range("D5").Value = 26.252149647524135
$ws.Range("B7").Value = 166.68748162564415
$ws.Range("C7").Value = 3508.3869994848328
$ws.Range("D7").Value = 1651.9886781408256
$ws.Range("B8").Value = 303.85093070767886
$ws.Range("C8").Value = 2515.7976568345643
$ws.Range("D8").Value = 575.40853887706555
$ws.Range("B11").Value = 81.841286077734097
$ws.Range("D11").Value = 270.11511184407948
$ws.Range("B12").Value = 1681.8239427398062
$ws.Range("C12").Value = 9085.581904209972
$ws.Range("D12").Value = 1662.9494410809871
$ws.Range("B14").Value = 1656.6838907742795
$ws.Range("C14").Value = 4134.7923456400731
$ws.Range("D14").Value = 625.63413573865967
$ws.Range("B17").Value = 17.192866011324565
$ws.Range("C17").Value = 491.09467172678018
$ws.Range("D17").Value = 471.41093256912598
$ws.Range("B18").Value = 136.23859446462069
$ws.Range("C18").Value = 2388.6728025644097
$ws.Range("D18").Value = 3234.0372034154675
$ws.Range("B19").Value = 80.199673262420561
$ws.Range("C19").Value = 2153.000516323043
$ws.Range("D19").Value = 1749.4869938197646
$ws.Range("B20").Value = 8363.7968380756793
$ws.Range("C20").Value = 30849.821655784439
$ws.Range("D20").Value = 4541.8512553986975
$ws.Range("B21").Value = 1285.6987388156779
$ws.Range("C21").Value = 6067.9139580607552
$ws.Range("D21").Value = 253.7027592337414
$ws.Range("B22").Value = 340.7977514284745
$ws.Range("C22").Value = 1235.1877670486772
$ws.Range("D22").Value = 408.6479771200701
$ws.Range("B23").Value = 548.27389376515089
$ws.Range("C23").Value = 6381.8123306183506
$ws.Range("D23").Value = 1530.284089409038
$ws.Range("B24").Value = 68.919612194910229
$ws.Range("C24").Value = 2904.849567855842
$ws.Range("D24").Value = 65.110955780769785
$ws.Range("B25").Value = 911.48137746857958
$ws.Range("C25").Value = 784.21464587640139
$ws.Range("D25").Value = 2432.8230458944254
$ws.Range("B26").Value = 1209.7157595749534
$ws.Range("C26").Value = 8895.7260631063164
$ws.Range("D26").Value = 659.39702925757194
$ws.Range("B27").Value = 6408.7615134398275
$ws.Range("C27").Value = 5627.2125656366197
$ws.Range("D27").Value = 1170.1316528285643
$ws.Range("B28").Value = 1782.3988299362693
$ws.Range("C28").Value = 7202.7636055719831
$ws.Range("D28").Value = 437.92287246104371

Write-Output "Updated emission factor values for data_egrid_emf_Cl sheet (paper resubmission revision)"
